# Updated cryptos list on Sun Nov 24 18:51:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.357.41"
$ws.Range("E2").Value = "  -1.57%  "

$ws.Range("D3").Value = "3.330.61"
$ws.Range("E3").Value = "  -2.83%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "249.77"
$ws.Range("E5").Value = "  -2.55%  "

$ws.Range("D6").Value = "'654.60"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "1.39"
$ws.Range("E7").Value = "  -6.96%  "

$ws.Range("E8").Value = "  -2.37%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -6.23%  "

$ws.Range("D11").Value = "3.329.23"
$ws.Range("E11").Value = "  -2.78%  "

$ws.Range("E12").Value = "  -3.66%  "

$ws.Range("D13").Value = "40.19"
$ws.Range("E13").Value = "  -4.52%  "

$ws.Range("D14").Value = "96.081.37"
$ws.Range("E14").Value = "  -1.59%  "

$ws.Range("E15").Value = "  -4.73%  "

$ws.Range("E16").Value = "  -3.60%  "

$ws.Range("D17").Value = "3.958.16"
$ws.Range("E17").Value = "  -2.60%  "

$ws.Range("E18").Value = "  -3.47%  "

$ws.Range("D19").Value = "3.339.09"
$ws.Range("E19").Value = "  -2.26%  "

$ws.Range("E20").Value = "  -3.05%  "

$ws.Range("D21").Value = "0.514"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D24").Value = "3.35"
$ws.Range("E24").Value = "  -3.28%  "

$ws.Range("E25").Value = "  -4.39%  "

$ws.Range("D26").Value = "6.55"
$ws.Range("E26").Value = "  +5.80%  "

$ws.Range("D27").Value = "96.05"
$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("E28").Value = "  -5.98%  "

$ws.Range("D29").Value = "3.504.77"
$ws.Range("E29").Value = "  -2.05%  "

$ws.Range("E30").Value = "  -8.76%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").Value = "11.07"
$ws.Range("E32").Value = "  -3.56%  "

$ws.Range("D33").Value = "0.187"
$ws.Range("E33").Value = "  -6.83%  "

$ws.Range("D34").Value = "2.48"
$ws.Range("E34").Value = "  +9.08%  "

$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("D36").Value = "0.546"
$ws.Range("E36").Value = "  -5.22%  "

$ws.Range("D37").Value = "'28.00"
$ws.Range("E37").Value = "  -7.17%  "

$ws.Range("E38").Value = "  +3.77%  "

$ws.Range("D39").Value = "7.61"

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("E41").Value = "  -3.78%  "

$ws.Range("D42").Value = "509.01"
$ws.Range("E42").Value = "  -1.81%  "

$ws.Range("D43").Value = "24.34"

$ws.Range("D44").Value = "0.831"
$ws.Range("E44").Value = "  -4.34%  "

$ws.Range("D47").Value = "1.66"
$ws.Range("E47").Value = "  +5.24%  "

$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("E49").Value = "  +1.57%  "

$ws.Range("D50").Value = "53.12"
$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("D51").Value = "3.12"
$ws.Range("E51").Value = "  -5.71%  "

# Row 22/23 swap: BitcoinCash <-> Uniswap
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "10.52"
$ws.Range("E22").Value = "  -4.91%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "503.41"
$ws.Range("E23").Value = "  -1.21%  "

# Row 45/46 swap: MantraDAO <-> VeChain
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0417"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").Value = "3.64"
$ws.Range("E46").Value = "  -1.22%  "
